# Remove pywin32 from requirements to deploy
# (workbook edit: append a duplicated header + summary block below the
#  existing results table, border the data rows, and tweak the view/print
#  setup to match)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append a second copy of the header row (row 26) and the last six
#     data rows (rows 20:25 -> 27:32) below the existing table. Using
#     Range.Copy(Destination) brings the values AND the existing cell
#     style along with it (reusing the header's bold/border/centered
#     style and leaving the freshly-appended data rows unstyled, same as
#     the source rows before the border pass below).
$ws.Range("A1:G1").Copy($ws.Range("A26"))
$ws.Range("A20:G25").Copy($ws.Range("A27"))

# --- Give every data row (old + newly appended) a thin box border. This
#     mirrors the new cellXfs entry (border-only style) being stamped
#     across A2:G25 and A27:G32 while leaving the header rows (1 and 26)
#     on their existing bold/centered style.
$ws.Range("A2:G25").Borders.LineStyle = 1
$ws.Range("A27:G32").Borders.LineStyle = 1

# --- Move the selection down to the newly appended summary block.
$ws.Range("A27:G32").Select()

# --- Set up the print page (paper size + portrait orientation).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
